$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value (tab-separated).
# A leading apostrophe marker is used so Excel stores the value as
# literal text even when it looks like a number (e.g. "103.20", "0.0978"),
# preserving exact formatting instead of coercing it to a numeric value.
$updates = @"
D2	42.604.94
E2	  +2.63%  
D3	2.197.70
E3	  +1.59%  
E4	  -0.03%  
D5	251.03
E5	  +5.71%  
E6	  +1.20%  
D7	74.62
E7	  +4.72%  
E8	  -0.08%  
D9	0.589
E9	  +2.54%  
E10	  +1.66%  
E11	  +2.03%  
B12	Polkadot
C12	https://coinranking.com/coin/25W7FG7om+polkadot-dot
D12	6.84
E12	  +2.46%  
B13	TRON
C13	https://coinranking.com/coin/qUhEFk1I61atv+tron-trx
D13	0.101
E13	  +1.61%  
D14	2.529.52
E14	  +1.65%  
D15	14.33
E15	  +1.40%  
D16	2.192.31
E16	  +1.63%  
D17	0.780
E17	  +0.11%  
D18	42.543.84
E18	  +2.86%  
E19	  +1.46%  
E20	  +2.30%  
D21	5.93
E21	  +3.03%  
D22	228.67
E22	  +0.60%  
E23	  +8.96%  
D24	9.45
E24	  -5.05%  
E25	  -0.19%  
D26	10.70
E26	  +0.19%  
D27	3.37
E27	  +3.66%  
E28	  +1.17%  
B29	InjectiveProtocol
C29	https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj
D29	37.48
E29	  +12.89%  
B30	Toncoin
C30	https://coinranking.com/coin/67YlI0K1b+toncoin-ton
D30	2.14
E30	  -1.08%  
D31	169.34
E31	  -1.31%  
D32	20.09
E32	  +1.82%  
E33	  +3.79%  
D34	5.17
E34	  +1.39%  
E35	  +1.07%  
E36	  +3.33%  
E37	  +3.88%  
E38	  +8.97%  
D39	12.20
E39	  +1.84%  
E40	  +0.76%  
E41	  +5.53%  
E42	  -1.00%  
D43	59.21
E43	  +0.93%  
D44	103.20
D45	0.481
E45	  +24.42%  
E46	  +1.28%  
D47	0.0978
D48	2.43
E48	  +12.98%  
E49	  +2.56%  
E50	  +1.96%  
E51	  +1.54%  
"@

foreach ($line in ($updates -split "`r?`n")) {
    if ([string]::IsNullOrEmpty($line)) { continue }
    $parts = $line -split "`t", 2
    $cellRef = $parts[0]
    $newValue = $parts[1]
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $newValue
    $range.Style = "Normal"
}
